# Adds the "countries_populations" CREATE TABLE ... AS (SELECT ... FROM
# countries INNER JOIN populations ON ...) statement after the existing
# "CREATE TABLE countries (...);" block at the end of the document.
#
# The new content is inserted, as raw WordprocessingML, immediately before
# the closing ");" of the last paragraph in the document. Range.InsertXML
# splits on the embedded <w:p> boundaries: every complete <w:p>...</w:p> in
# the payload becomes its own new paragraph, while the payload's own final
# (unclosed-relative-to-destination) fragment is merged into the start of
# the destination paragraph -- so the original ");" paragraph survives
# untouched (same paraId/rsids) with the new "ON ... = ...country_id"
# runs simply prepended to its existing ");" run.

$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$wNs = "xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`""

$newContentXml = (
    "<w:p $wNs><w:r><w:t>);</w:t></w:r></w:p>" +
    "<w:p $wNs/>" +
    "<w:p $wNs/>" +
    "<w:p $wNs>" +
        "<w:r><w:t xml:space=`"preserve`">CREATE TABLE </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:t>countries_populations</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`"> AS</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wNs>" +
        "<w:r><w:tab/><w:t xml:space=`"preserve`">(SELECT </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:proofErr w:type=`"gramStart`"/>" +
        "<w:r><w:t>populations.record</w:t></w:r>" +
        "<w:proofErr w:type=`"gramEnd`"/>" +
        "<w:r><w:t>_key</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`">, </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:t>populations.country_id</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`">, </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:t>countries.country_name</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`">, </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:t>populations.mid_year</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`">, </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:t>populations.population</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
    "</w:p>" +
    "<w:p $wNs>" +
        "<w:r><w:tab/><w:t>FROM countries</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wNs>" +
        "<w:r><w:tab/></w:r>" +
        "<w:r><w:tab/><w:t xml:space=`"preserve`">INNER JOIN populations </w:t></w:r>" +
    "</w:p>" +
    "<w:p $wNs>" +
        "<w:r><w:tab/><w:t xml:space=`"preserve`"> </w:t></w:r>" +
        "<w:r><w:tab/></w:r>" +
        "<w:r><w:tab/><w:t xml:space=`"preserve`">ON </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:t>countries.country_id</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`"> = </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:t>populations.country_id</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
    "</w:p>"
)

$insertionPoint.InsertXML($newContentXml)
